$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "得分" (Score) header to "平均分" (Average score)
$ws.Range("E1").Value = "平均分"

# Re-sort the city data rows (2-28) in descending order of the average
# score column (E), keeping each city's A:E row intact.
$dataRange = $ws.Range("A2:E28")
$sortKey = $ws.Range("E2:E28")
$dataRange.Sort($sortKey, 2)
